$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$rsquo = [char]0x2019

# --- Shape id=142 "Formulation of new ecPoint-Rainfall products..." ---
# Split the sentence into two runs: rewrite the trailing clause and
# append the new closing sentence as a separate run.
$sp142 = Get-ShapeById $s 142
$tr142 = $sp142.TextFrame.TextRange
$tr142.Text = "Formulation of new ecPoint-Rainfall products and guidelines tailored to participants" + $rsquo + " needs discussed at "
$tr142.InsertAfter("previous step.") | Out-Null

# Shape id=142 line color FF00FF -> 00CCFF
$sp142.Line.ForeColor.RGB = 16763904

# --- Shape id=143 "Joint revision..." -> "Joint review..." ---
$sp143 = Get-ShapeById $s 143
$sp143.TextFrame.TextRange.Text = "Joint review of the summary reports results under the revised ecPoint-Rainfall products and guidelines."

# --- Shape id=156 "Independent revision..." -> "Independent review..." ---
$sp156 = Get-ShapeById $s 156
$sp156.TextFrame.TextRange.Text = "Independent review of the summary reports" + $rsquo + " content."
